# "Work on the Arduino Watch and Final Report"
#
# Applies:
#  - New Print_Area defined name (Sheet1!$AH$5:$AJ$15) and matching
#    selection/active-cell on the sheet (AH5, sqref AH5:AJ15)
#  - Re-colours a handful of Gantt-bar cells:
#      AA8  (Database week4) : green  -> default   (task no longer runs that week)
#      AA10 (Servers  week4) : blue   -> default   (task no longer runs that week)
#      Y12  (row12  week3)   : default-> red/delay
#      Z12  (row12  week4)   : default-> red/delay
#      Y16  (row16  week3)   : orange -> red/delay

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set the print area (creates the workbook-level defined name) ---
$ws.PageSetup.PrintArea = '$AH$5:$AJ$15'

# --- Move the visible selection to the new print area ---
$ws.Range("AH5:AJ15").Select()

# --- Re-colour Gantt cells by copying formats from donor cells that ---
# --- already carry the desired style, so the existing style records ---
# --- are reused instead of new (duplicate) ones being minted.       ---

# AA8 & AA10 revert to the plain/default cell background (same as AA9, AA11, ...)
$ws.Range("AA9").Copy()
$ws.Range("AA8").PasteSpecial(-4122)

$ws.Range("AA11").Copy()
$ws.Range("AA10").PasteSpecial(-4122)

# Y12, Z12 and Y16 take on the "delay" colour already used at Q4
$ws.Range("Q4").Copy()
$ws.Range("Y12").PasteSpecial(-4122)
$ws.Range("Q4").Copy()
$ws.Range("Z12").PasteSpecial(-4122)
$ws.Range("Q4").Copy()
$ws.Range("Y16").PasteSpecial(-4122)

$excel.CutCopyMode = 0
